$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-12 19:18:50'
$ws.Range('E3').Value = '2026-02-12 19:18:52'
$ws.Range('E4').Value = '2026-02-12 19:18:55'
$ws.Range('H4').Value = '''38%'
$ws.Range('J4').Value = '998.8 hPa'
$ws.Range('K4').Value = '13.3 MJ/m2'
$ws.Range('O4').Value = '16.3 °C'
$ws.Range('E5').Value = '2026-02-12 19:18:58'
$ws.Range('O5').Value = '-3.3 °C'
$ws.Range('E6').Value = '2026-02-12 19:19:01'
$ws.Range('J6').Value = '998.6 hPa'
$ws.Range('E7').Value = '2026-02-12 19:19:03'
$ws.Range('J7').Value = '1001.5 hPa'
$ws.Range('E8').Value = '2026-02-12 19:19:06'
$ws.Range('H8').Value = '''45%'
$ws.Range('J8').Value = '1000.8 hPa'
$ws.Range('E9').Value = '2026-02-12 19:19:09'
$ws.Range('E10').Value = '2026-02-12 19:19:11'
$ws.Range('E11').Value = '2026-02-12 19:19:14'
$ws.Range('H11').Value = '''47%'
$ws.Range('E12').Value = '2026-02-12 19:19:17'
$ws.Range('H12').Value = '''73%'
$ws.Range('E13').Value = '2026-02-12 19:19:19'
$ws.Range('J13').Value = '1001.4 hPa'
$ws.Range('E14').Value = '2026-02-12 19:19:22'
$ws.Range('E15').Value = '2026-02-12 19:19:25'
$ws.Range('H15').Value = '''55%'
$ws.Range('O15').Value = '14.0 °C'
$ws.Range('E16').Value = '2026-02-12 19:19:27'
$ws.Range('E17').Value = '2026-02-12 19:19:30'
$ws.Range('O17').Value = '2.1 °C'
$ws.Range('E18').Value = '2026-02-12 19:19:33'
$ws.Range('J18').Value = '999.1 hPa'
$ws.Range('E19').Value = '2026-02-12 19:19:35'
$ws.Range('E20').Value = '2026-02-12 19:19:38'
$ws.Range('E21').Value = '2026-02-12 19:19:41'
$ws.Range('J21').Value = '1001.8 hPa'
$ws.Range('E22').Value = '2026-02-12 19:19:43'
$ws.Range('O22').Value = '-5.6 °C'
$ws.Range('E23').Value = '2026-02-12 19:19:46'
$ws.Range('E24').Value = '2026-02-12 19:19:49'
$ws.Range('E25').Value = '2026-02-12 19:19:52'
$ws.Range('H25').Value = '''60%'
$ws.Range('O25').Value = '-1.9 °C'
$ws.Range('E26').Value = '2026-02-12 19:19:54'
$ws.Range('J26').Value = '998.1 hPa'
$ws.Range('O26').Value = '6.1 °C'
$ws.Range('E27').Value = '2026-02-12 19:19:57'
$ws.Range('E28').Value = '2026-02-12 19:20:00'
$ws.Range('J28').Value = '998.4 hPa'
$ws.Range('E29').Value = '2026-02-12 19:20:02'
$ws.Range('N29').Value = '10.2 °C 18:39 TU'
$ws.Range('O29').Value = '15.1 °C'
$ws.Range('E30').Value = '2026-02-12 19:20:05'
$ws.Range('J30').Value = '998.9 hPa'
$ws.Range('O30').Value = '12.5 °C'
$ws.Range('E31').Value = '2026-02-12 19:20:08'
$ws.Range('J31').Value = '998.3 hPa'
$ws.Range('O31').Value = '14.4 °C'
$ws.Range('E32').Value = '2026-02-12 19:20:10'
$ws.Range('E33').Value = '2026-02-12 19:20:13'
$ws.Range('J33').Value = '1001.0 hPa'
$ws.Range('E34').Value = '2026-02-12 19:20:16'
$ws.Range('E35').Value = '2026-02-12 19:20:18'
$ws.Range('J35').Value = '1007.8 hPa'
$ws.Range('E36').Value = '2026-02-12 19:20:21'
$ws.Range('H36').Value = '''62%'
$ws.Range('J36').Value = '999.2 hPa'
$ws.Range('E37').Value = '2026-02-12 19:20:23'
$ws.Range('H37').Value = '''47%'
$ws.Range('J37').Value = '999.7 hPa'
$ws.Range('O37').Value = '10.3 °C'
$ws.Range('E38').Value = '2026-02-12 19:20:26'
$ws.Range('E39').Value = '2026-02-12 19:20:29'
$ws.Range('E40').Value = '2026-02-12 19:20:32'
$ws.Range('J40').Value = '1002.6 hPa'
$ws.Range('O40').Value = '9.9 °C'
$ws.Range('E41').Value = '2026-02-12 19:20:34'
$ws.Range('J41').Value = '1005.4 hPa'
$ws.Range('K41').Value = '14.1 MJ/m2'
$ws.Range('E42').Value = '2026-02-12 19:20:37'
$ws.Range('O42').Value = '14.5 °C'
$ws.Range('E43').Value = '2026-02-12 19:20:39'
$ws.Range('E44').Value = '2026-02-12 19:20:42'
$ws.Range('E45').Value = '2026-02-12 19:20:45'
$ws.Range('J45').Value = '1004.6 hPa'
$ws.Range('O45').Value = '7.3 °C'
$ws.Range('E46').Value = '2026-02-12 19:20:48'
$ws.Range('N46').Value = '11.8 °C 18:51 TU'
$ws.Range('O46').Value = '16.1 °C'
